# Notion "Chấm công HỆ THỐNG" export - apply the latest Notion sync update.
# - last_edited_time for rows 2-13 moves to 2024-08-03T20:14:00.000Z
# - last_edited_time for rows 14-22 moves to 2024-08-03T20:15:00.000Z
# - A handful of computed "công" totals (Tổng công tại LONG XUYÊN / Nửa ngày /
#   Tổng công tại CẦN THƠ / Tổng công) were recalculated after the multi
#   process strategy change, so their numeric values change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last_edited_time (column D) -------------------------------------------------
$lastEditedTime1 = "2024-08-03T20:14:00.000Z"
$lastEditedTime2 = "2024-08-03T20:15:00.000Z"

$ws.Range("D2:D13").Value = $lastEditedTime1
$ws.Range("D14:D22").Value = $lastEditedTime2

# --- recalculated totals ----------------------------------------------------------
$ws.Range("S10").Value = 2.5
$ws.Range("V10").Value = 1
$ws.Range("AM10").Value = 2.5

$ws.Range("S16").Value = 3
$ws.Range("AF16").Value = 3
$ws.Range("AM16").Value = 3

$ws.Range("AI18").Value = 0.5
$ws.Range("AM18").Value = 0.5

$ws.Range("S19").Value = 3
$ws.Range("AF19").Value = 3
$ws.Range("AM19").Value = 3

$ws.Range("S21").Value = 3
$ws.Range("AF21").Value = 3
$ws.Range("AM21").Value = 3

$ws.Range("S22").Value = 3
$ws.Range("AF22").Value = 3
$ws.Range("AM22").Value = 3

$wb.Save()
